$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add a date value to C1 (serial 44307 = 2021-04-21), formatted as a date
# (numFmtId 14 - matches the new cellXf introduced in styles.xml).
$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
